$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Step 0: remove the existing _GoBack bookmark (it will be re-inserted
# later into its own empty paragraph).
# -----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -----------------------------------------------------------------------
# Step 1: split paragraph 4 ("https://expressjs.com") to create two new
# empty ListParagraph paragraphs *before* it:
#   para4 (new) -> will hold the expressjs.com hyperlink
#   para5 (new) -> will hold the re-inserted _GoBack bookmark
# para6 keeps the original "https://expressjs.com" text for now, which
# we change afterwards.
# -----------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$insPoint = $d.Range($p4.Range.Start, $p4.Range.Start)
$insPoint.InsertBefore("`r`r")

# Re-resolve paragraphs after the insert.
$linkPara1 = $d.Paragraphs(4)
$bmPara = $d.Paragraphs(5)
$textPara = $d.Paragraphs(6)

# Insert the expressjs.com/ hyperlink into the (currently empty) paragraph 4.
$d.Hyperlinks.Add($linkPara1.Range, "https://expressjs.com/", $null, $null, "https://expressjs.com/") | Out-Null

# Re-insert the _GoBack bookmark, collapsed, inside the empty paragraph 5.
$bmRange = $d.Range($bmPara.Range.Start, $bmPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Replace the text of the old "https://expressjs.com" paragraph.
$textPara.Range.Text = "UI UX giao dien moble"

# -----------------------------------------------------------------------
# Step 2: after the "UI UX giao dien moble" paragraph, append:
#   - a paragraph with the material-ui.com hyperlink
#   - three empty ListParagraph paragraphs
# -----------------------------------------------------------------------
$textPara = $d.Paragraphs(6)
$endPoint = $d.Range($textPara.Range.End - 1, $textPara.Range.End - 1)
$endPoint.InsertAfter("`r`r`r`r")

$linkPara2 = $d.Paragraphs(7)
$d.Hyperlinks.Add($linkPara2.Range, "https://material-ui.com/styles/advanced/", $null, $null, "https://material-ui.com/styles/advanced/") | Out-Null

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    Write-Host "Para $i [$($p.Range.Text)]"
}
